# Monthly task assignment workbook update
# - Pipette/volumetric task values now carry unit suffixes (mL) and amber/clear labels
# - "Analyst N" section headers get a slightly larger bold font (11 -> 12 pt)
# - Columns A, B and D get explicit widths so the longer values are readable

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Pipette task values (column B) and Volumetric task values (column D) ---

$ws.Range("B4").Value = "5mL"
$ws.Range("D4").Value = "100mL amber"
$ws.Range("B5").Value = "2mL"
$ws.Range("D5").Value = "25mL clear"
$ws.Range("B6").Value = "1mL"
$ws.Range("D6").Value = "250mL clear"
$ws.Range("B7").Value = "10mL"
$ws.Range("D7").Value = "100mL amber"
$ws.Range("B8").Value = "4mL"
$ws.Range("D8").Value = "50mL clear"
$ws.Range("B9").Value = "2mL"
$ws.Range("D9").Value = "100mL amber"

$ws.Range("B12").Value = "1mL"
$ws.Range("D12").Value = "50mL amber"
$ws.Range("B13").Value = "5mL"
$ws.Range("D13").Value = "250mL clear"
$ws.Range("B14").Value = "2mL"
$ws.Range("D14").Value = "25mL clear"
$ws.Range("B15").Value = "20mL"
$ws.Range("D15").Value = "100mL clear"
$ws.Range("B16").Value = "2mL"
$ws.Range("D16").Value = "25mL clear"
$ws.Range("B17").Value = "1mL"
$ws.Range("D17").Value = "50mL amber"

$ws.Range("B20").Value = "20mL"
$ws.Range("D20").Value = "25mL clear"
$ws.Range("B21").Value = "10mL"
$ws.Range("D21").Value = "100mL amber"
$ws.Range("B22").Value = "4mL"
$ws.Range("D22").Value = "50mL clear"
$ws.Range("B23").Value = "5mL"
$ws.Range("D23").Value = "25mL clear"
$ws.Range("B24").Value = "20mL"
$ws.Range("D24").Value = "50mL amber"
$ws.Range("B25").Value = "4mL"
$ws.Range("D25").Value = "100mL clear"

$ws.Range("B28").Value = "2mL"
$ws.Range("D28").Value = "50mL clear"
$ws.Range("B29").Value = "20mL"
$ws.Range("D29").Value = "100mL clear"
$ws.Range("B30").Value = "10mL"
$ws.Range("D30").Value = "50mL amber"
$ws.Range("B31").Value = "2mL"
$ws.Range("D31").Value = "50mL clear"
$ws.Range("B32").Value = "10mL"
$ws.Range("D32").Value = "100mL clear"
$ws.Range("B33").Value = "20mL"
$ws.Range("D33").Value = "25mL clear"

$ws.Range("B36").Value = "10mL"
$ws.Range("D36").Value = "100mL clear"
$ws.Range("B37").Value = "4mL"
$ws.Range("D37").Value = "50mL amber"
$ws.Range("B38").Value = "5mL"
$ws.Range("D38").Value = "100mL amber"
$ws.Range("B39").Value = "1mL"
$ws.Range("D39").Value = "250mL clear"
$ws.Range("B40").Value = "5mL"
$ws.Range("D40").Value = "100mL amber"
$ws.Range("B41").Value = "10mL"
$ws.Range("D41").Value = "250mL clear"

$ws.Range("B44").Value = "4mL"
$ws.Range("D44").Value = "250mL clear"
$ws.Range("B45").Value = "1mL"
$ws.Range("D45").Value = "50mL clear"
$ws.Range("B46").Value = "20mL"
$ws.Range("D46").Value = "100mL clear"
$ws.Range("B47").Value = "4mL"
$ws.Range("D47").Value = "50mL amber"
$ws.Range("B48").Value = "1mL"
$ws.Range("D48").Value = "250mL clear"
$ws.Range("B49").Value = "5mL"
$ws.Range("D49").Value = "50mL clear"

# --- Bold "Analyst N" header cells now use a 12pt font instead of 11pt ---

$ws.Range("B3").Font.Size = 12
$ws.Range("B11").Font.Size = 12
$ws.Range("B19").Font.Size = 12
$ws.Range("B27").Font.Size = 12
$ws.Range("B35").Font.Size = 12
$ws.Range("B43").Font.Size = 12

# --- Column widths so the new, longer values fit comfortably ---

$ws.Columns("A:B").ColumnWidth = 18
$ws.Columns("D").ColumnWidth = 22
